$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 14: "SRv6 with Timestamp and Forward Function"
#   - update "Protocol = UDP" line to "Next Header = 43 (Routing Header)"
#   - add two new lines ("Next Header = 17 (UDP)" + blank) before the
#     next "+---+" separator
#   - update the figure caption
#   - the diagram textbox ("Rectangle 5") auto-fits to the extra lines;
#     the height is (re)applied explicitly afterwards to match exactly
# ---------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$rect14 = $s14.Shapes.Item("Rectangle 5")
$tr14 = $rect14.TextFrame.TextRange

$protoPara = $tr14.Paragraphs(5,1)
$protoRun = $protoPara.Runs(1,1)
$protoRun.Text = "  .  Next Header = 43 (Routing Header)                            ."

$blankPara = $tr14.Paragraphs(6,1)
$null = $blankPara.InsertAfter("`r  .  Next Header = 17 (UDP)                                       .`r  .                                                               .")

$figPara14 = $tr14.Paragraphs(24,1)
$figRun14 = $figPara14.Runs(1,1)
$figRun14.Text = "    Figure: Example Probe Message for SRv6 with Endpoint Function"

# Resize the box (off stays the same, only the height grows to fit the
# two extra lines of text)
$rect14.Height = (3785652 / 12700) + 0.00001

# ---------------------------------------------------------------------
# Slide 8: "SR-MPLS with Timestamp Label"
#   - update the figure caption
#   - reposition/resize the diagram textbox ("Rectangle 5")
#   - reposition the adjacent bullet textbox ("Content Placeholder 7")
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$rect8 = $s8.Shapes.Item("Rectangle 5")
$tr8 = $rect8.TextFrame.TextRange

$figPara8 = $tr8.Paragraphs(29,1)
$figRun8 = $figPara8.Runs(1,1)
$figRun8.Text = "   Figure 5: Example Probe Message for SR-MPLS with Timestamp Label"

$rect8.Left = (173421 / 12700) + 0.00001
$rect8.Top = (650498 / 12700) + 0.00001
$rect8.Width = (5105400 / 12700) + 0.00001
$rect8.Height = (4247317 / 12700) + 0.00001

$content8 = $s8.Shapes.Item("Content Placeholder 7")
$content8.Left = (5278821 / 12700) + 0.00001
$content8.Top = (1071222 / 12700) + 0.00001
